# Add a new "Note:publication" column (R) to the "Item description" sheet,
# immediately to the right of the existing "Note:preferred citation" column (Q),
# and record the related publication citation for the second data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Item description")
$ws.Activate()

# Insert a new column at R; this shifts R:AG -> S:AH and copies the
# formatting of the column immediately to the left (Q) onto the new cells.
$ws.Columns("R:R").Insert()

# Match the new column's width to its left neighbour (Q), same as the other
# "Note:*" columns in this block.
$ws.Columns("R:R").ColumnWidth = $ws.Columns("Q:Q").ColumnWidth

# Header for the new column.
$ws.Range("R1").Value = "Note:publication"

# Related-publication citation for row 2 (Revillagigedo / rocky-reef health index item).
$ws.Range("R2").Value = "Aburto-Oropeza O, et al (2014) A framework to assess the health of rocky reefs linking geomorphology, community assemblage, and fish biomass. Ecological Indicators 52:353-361. http://dx.doi.org/10.1016/j.ecolind.2014.12.006."

# Row 3 has no publication note (left blank), matching the other empty cells
# in that column style.

# Update the visible selection to reflect where editing left off.
$ws.Range("Q2").Select()
